$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared text updates (rich-text runs) ---
# A8 / "Volume 31   Number  25" -> "...26"
$ws.Range("A8").Characters(21,2).Text = "26"
# C9 / "Report Covering the Week  6/17/2024  Through  6/23/2024"
$ws.Range("C9").Characters(27,9).Text = "6/24/2024"
$ws.Range("C9").Characters(47,9).Text = "6/30/2024"

# --- Reference cells used as style/shared-string donors (never modified themselves) ---
# C14 = shared "0" marker (style 14); M14 = shared "***.*" marker (style 14)
# G14 = plain integer style 15; H14 = plain percent style 16

# --- Row 14 ---
$ws.Range("C14").Copy($ws.Range("D14"))   # D14: 1 -> "0" marker
$ws.Range("M14").Copy($ws.Range("E14"))   # E14: -100 -> "***.*" marker

# --- Row 16 ---
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 42.857142857142
$ws.Range("I16").Value = 55
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 37.5
$ws.Range("L16").Value = 37.5
$ws.Range("M16").Value = -32.098765432098
$ws.Range("N16").Value = -84.931506849315

# --- Row 17 ---
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = 16.666666666666
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -30.434782608695
$ws.Range("I17").Value = 108
$ws.Range("J17").Value = 109
$ws.Range("K17").Value = -0.917431192660
$ws.Range("L17").Value = 3.846153846153
$ws.Range("M17").Value = 61.194029850746
$ws.Range("N17").Value = -29.870129870129

# --- Row 18 ---
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -36.363636363636
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 59
$ws.Range("K18").Value = -5.084745762711
$ws.Range("L18").Value = -32.530120481927
$ws.Range("M18").Value = -70.526315789473
$ws.Range("N18").Value = -93.128834355828

# --- Row 19 ---
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 15
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -10.416666666666
$ws.Range("I19").Value = 270
$ws.Range("J19").Value = 280
$ws.Range("K19").Value = -3.571428571428
$ws.Range("L19").Value = -8.474576271186
$ws.Range("M19").Value = 39.896373056994
$ws.Range("N19").Value = -10.891089108910

# --- Row 20 ---
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 40
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -6.25
$ws.Range("I20").Value = 89
$ws.Range("J20").Value = 76
$ws.Range("K20").Value = 17.105263157894
$ws.Range("L20").Value = 85.416666666666
$ws.Range("M20").Value = 23.611111111111
$ws.Range("N20").Value = -90.541976620616

# --- Row 21 ---
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 106
$ws.Range("H21").Value = -14.150943396226
$ws.Range("I21").Value = 588
$ws.Range("J21").Value = 576
$ws.Range("K21").Value = 2.083333333333
$ws.Range("L21").Value = 1.030927835051
$ws.Range("M21").Value = -2.809917355371
$ws.Range("N21").Value = -77.323563440030

# --- Row 22 ---
$ws.Range("G14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = 16.666666666666
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = 0

# --- Row 24 ---
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -45.833333333333
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = -17.204301075268
$ws.Range("I24").Value = 494
$ws.Range("J24").Value = 568
$ws.Range("K24").Value = -13.028169014084
$ws.Range("L24").Value = -10.669077757685
$ws.Range("M24").Value = 11.764705882352

# --- Row 25 ---
$ws.Range("C14").Copy($ws.Range("C25"))
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 14.285714285714
$ws.Range("J25").Value = 160
$ws.Range("K25").Value = -25.625
$ws.Range("L25").Value = -28.313253012048

# --- Row 26 ---
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 34
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = 6.25
$ws.Range("I26").Value = 204
$ws.Range("J26").Value = 168
$ws.Range("K26").Value = 21.428571428571
$ws.Range("L26").Value = 27.5
$ws.Range("M26").Value = -10.132158590308

# --- Row 28 ---
$ws.Range("C28").Value = 3
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 37
$ws.Range("J28").Value = 32
$ws.Range("K28").Value = 15.625
$ws.Range("L28").Value = 0

# --- Row 29 ---
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("M14").Copy($ws.Range("E29"))

# --- Row 30 ---
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("M14").Copy($ws.Range("E30"))

# --- Row 31 ---
$ws.Range("G14").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("H14").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("G14").Copy($ws.Range("G31"))
$ws.Range("G31").Value = 1
$ws.Range("H14").Copy($ws.Range("H31"))
$ws.Range("H31").Value = 100
$ws.Range("J31").Value = 4
$ws.Range("K31").Value = 150

# --- Row 33 ---
$ws.Range("G14").Copy($ws.Range("D33"))
$ws.Range("D33").Value = 1
$ws.Range("H14").Copy($ws.Range("E33"))
$ws.Range("E33").Value = -100
$ws.Range("G33").Value = 2
$ws.Range("J33").Value = 4
$ws.Range("K33").Value = -75
$ws.Range("L33").Value = -50
